$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet/tab name (Excel enforces a 31-character limit on sheet
# names, so the full target name is truncated to fit that limit)
$ws.Name = " Leftover Pellets Over The Past"

# Clear old extra rows (9, 10, 11) data since new data only spans 8 rows
$ws.Rows.Item(9).ClearContents()
$ws.Rows.Item(10).ClearContents()
$ws.Rows.Item(11).ClearContents()

# Header row
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "first feed number of pellets left"
$ws.Range("C1").Value = "second feed number of pellets left"
$ws.Range("D1").Value = "total feed pellets fed"

# Copy style of A1/B1 (bold header) to C1/D1
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1:D1").PasteSpecial(-4122) | Out-Null

# Data rows
$ws.Range("A2").Value = "20 May"
$ws.Range("B2").Value = 25
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 400

$ws.Range("A3").Value = "21 May"
$ws.Range("B3").Value = 15
$ws.Range("C3").Value = 40
$ws.Range("D3").Value = 300

$ws.Range("A4").Value = "22 May"
$ws.Range("B4").Value = 25
$ws.Range("C4").Value = 30
$ws.Range("D4").Value = 250

$ws.Range("A5").Value = "23 May"
$ws.Range("B5").Value = 25
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 400

$ws.Range("A6").Value = "24 May"
$ws.Range("B6").Value = 15
$ws.Range("C6").Value = 40
$ws.Range("D6").Value = 300

$ws.Range("A7").Value = "25 May"
$ws.Range("B7").Value = 45
$ws.Range("C7").Value = 35
$ws.Range("D7").Value = 275

$ws.Range("A8").Value = "26 May"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
